$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# --- Step 1: insert the 5 new logbook paragraphs right after the "17/11/16" paragraph ---
$insertAfterPara = $d.Paragraphs.Item(6)
$insertAfterRange = $insertAfterPara.Range
$insertAfterRange.InsertParagraphAfter()

$target = $d.Paragraphs.Item(7).Range
$newParasXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>21/11/’16</w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t>de index pagina afwerken en de navigatie activeren</w:t></w:r></w:p><w:p><w:r><w:t>24/11/16</w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> navigatiebalk toevoegen en een aantal ervan activeren,</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">        Kleuren samenstellen en bepalen.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">        Homepage gemaakt in het Nederlands en in het frans.</w:t></w:r></w:p><w:p><w:r><w:t>29/11/16</w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> de breedte aangepast van mijn pagina indien deze word geopend op kleinere schermen</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($newParasXml)

# --- Step 2: split the bookmark paragraph so the " " run lives in its own new paragraph ---
# After step 1, the bookmark paragraph ("_GoBack") is paragraph 12.
$bmPara = $d.Paragraphs.Item(12)
$bmRange = $bmPara.Range
$bmRange.Collapse(0)
$bmRange.InsertParagraphAfter()

# Move the " " text into the freshly created (now empty) following paragraph.
$newPara = $d.Paragraphs.Item(13)
$newPara.Range.InsertAfter(" ")

# Remove the original " " run text that is still left at the end of the bookmark paragraph.
$bmPara2 = $d.Paragraphs.Item(12)
$delRange = $d.Range($bmPara2.Range.End - 2, $bmPara2.Range.End - 1)
$delRange.Delete()

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("Para " + $i + ": [" + $p.Range.Text + "]")
}
